$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" -> rename header B1 ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

# --- Sheet 2: "Monthly Trend" -> rename header B1 ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new sheet "PO Forecast" after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws3.Name = "PO Forecast"

# Match page margins of the other sheets (PageSetup uses points; XML stores inches)
$ws3.PageSetup.LeftMargin = 0.75 * 72
$ws3.PageSetup.RightMargin = 0.75 * 72
$ws3.PageSetup.TopMargin = 1 * 72
$ws3.PageSetup.BottomMargin = 1 * 72
$ws3.PageSetup.HeaderMargin = 0.5 * 72
$ws3.PageSetup.FooterMargin = 0.5 * 72

# Copy the header style (bold, border, centered) from an existing header cell
$ws1.Range("B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-column style (yyyy-mm-dd date number format) from an existing date cell
$ws1.Range("A2").Copy()
$ws3.Range("A2:A14").PasteSpecial(-4122)

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Data rows
$data = New-Object 'object[,]' 13,4
$data[0,0] = 45417.99999999999;  $data[0,1] = 18; $data[0,2] = 1.835149522489635;   $data[0,3] = 33.27024768481198
$data[1,0] = 45424.99999999999;  $data[1,1] = 16; $data[1,2] = 0.6240151608284173;  $data[1,3] = 32.27532836335961
$data[2,0] = 45431.99999999999;  $data[2,1] = 15; $data[2,2] = -1.815066650567576;  $data[2,3] = 32.14945974559848
$data[3,0] = 45438.99999999999;  $data[3,1] = 13; $data[3,2] = -2.475982074263318;  $data[3,3] = 30.65625377204262
$data[4,0] = 45494.99999999999;  $data[4,1] = 1;  $data[4,2] = -15.8038517956361;   $data[4,3] = 18.4587587481316
$data[5,0] = 45501.99999999999;  $data[5,1] = 0;  $data[5,2] = -16.1673933751304;   $data[5,3] = 16.46513037796724
$data[6,0] = 45508.99999999999;  $data[6,1] = 0;  $data[6,2] = -19.80625967068524;  $data[6,3] = 13.99831610427511
$data[7,0] = 45515.99999999999;  $data[7,1] = 0;  $data[7,2] = -20.28285281250611;  $data[7,3] = 11.75597852631919
$data[8,0] = 45522.99999999999;  $data[8,1] = 0;  $data[8,2] = -22.82321950552809;  $data[8,3] = 10.93223492568378
$data[9,0] = 45529.99999999999;  $data[9,1] = 0;  $data[9,2] = -22.80360858577489;  $data[9,3] = 10.3662793131795
$data[10,0] = 45536.99999999999; $data[10,1] = 0; $data[10,2] = -25.21547675962532; $data[10,3] = 8.299025849975017
$data[11,0] = 45543.99999999999; $data[11,1] = 0; $data[11,2] = -26.32504224274184; $data[11,3] = 7.251465701979724
$data[12,0] = 45550.99999999999; $data[12,1] = 0; $data[12,2] = -27.71743355393406; $data[12,3] = 4.846647857042933

$ws3.Range("A2:D14").Value = $data
